# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.283.96'
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').Value = '2.613.24'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'584.55"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.85%  '
$ws.Range('D6').Value = "'143.29"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.76%  '
$ws.Range('D7').Value = "'0.998"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').Value = "'0.598"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('D9').Value = "'6.50"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('E11').Value = '  +2.04%  '
$ws.Range('E12').Value = '  +0.96%  '
$ws.Range('D13').Value = '3.074.07'
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('D14').Value = "'25.20"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.39%  '
$ws.Range('D15').Value = '60.276.05'
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('D17').Value = '2.617.80'
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('E18').Value = '  +2.54%  '
$ws.Range('E19').Value = '  +0.91%  '
$ws.Range('D20').Value = "'346.51"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range('E21').Value = '  -2.09%  '
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('E23').Value = '  +1.17%  '
$ws.Range('D24').Value = "'63.67"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.89%  '
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('E26').Value = '  +0.32%  '
$ws.Range('D27').Value = "'8.04"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.51%  '
$ws.Range('D28').Value = "'1.95"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.92%  '
$ws.Range('E29').Value = '  +1.25%  '
$ws.Range('E30').Value = '  +2.32%  '
$ws.Range('E31').Value = '  +3.78%  '
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('D33').Value = "'19.51"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('E34').Value = '  +6.63%  '
$ws.Range('D35').Value = "'1.30"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.06%  '
$ws.Range('D36').Value = "'4.29"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.24%  '
$ws.Range('E37').Value = '  +2.24%  '
$ws.Range('E38').Value = '  +8.11%  '
$ws.Range('D39').Value = "'38.47"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.91%  '
$ws.Range('D40').Value = "'3.95"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.21%  '
$ws.Range('D41').Value = "'0.852"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('D42').Value = "'135.70"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.99%  '
$ws.Range('D43').Value = "'20.02"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.66%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = "'0.999"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.26%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').Value = "'0.0990"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.55%  '
$ws.Range('E46').Value = '  +0.82%  '
$ws.Range('E47').Value = '  +3.43%  '
$ws.Range('E48').Value = '  +1.29%  '
$ws.Range('E49').Value = '  +1.78%  '
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('D51').Value = "'10.75"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.44%  '
